# async.docx -- "Added fixes for Router 4.0"
#
# Two content edits:
#
#  1. "Install Async.js" section: the command the reader is told to type
#     changes from `npm install async` to `npm install --save async`
#     (typed/edited in separate bursts, so it ends up as three adjacent
#     InlineCode runs: "npm install" / " --save" / " async", with Word's
#     auto-managed "_GoBack" (last-edit-location) bookmark sitting
#     between the last two runs).
#
#  2. "Add the write and close tasks" section: the old edit-location
#     artifact (the "_GoBack" bookmark that used to sit here, splitting
#     "minification" from " workflow by adding...") is gone now that the
#     edit point has moved elsewhere, so those two runs collapse back
#     into a single run.
#
# Because "_GoBack" can only exist in one place at a time, moving it to
# its new home (edit 1) automatically removes it from its old home
# (edit 2).

$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Edit 2 first (while "_GoBack" still sits at its old location, use a
# temporary bookmark to protect the run boundary we want to *keep*,
# then force-merge the boundary we want to remove).
# -----------------------------------------------------------------

# Relocate "_GoBack" off to the document start for now -- this frees up
# its old slot (between "minification" and " workflow...") so those
# runs are eligible to recombine.
$d.Bookmarks.Add("_GoBack", $d.Range(0, 0))

$full = $d.Content.Text
$bText = "minification"
$idxB = $full.IndexOf("minification workflow by adding two more tasks")
$idxBC = $idxB + $bText.Length

# Protect the boundary before "minification" (between it and "We're now
# going to finish our ") so the upcoming merge doesn't swallow that run
# too.
$d.Bookmarks.Add("TEMP_PROTECT_A", $d.Range($idxB, $idxB))

# Insert then immediately delete a throw-away character right on the
# "minification"/" workflow..." boundary -- this round-trip is what
# triggers the two like-formatted runs either side of it to recombine
# into one.
$d.Range($idxBC, $idxBC).InsertAfter("x")
$d.Range($idxBC, $idxBC + 1).Delete()

$d.Bookmarks("TEMP_PROTECT_A").Delete()

# -----------------------------------------------------------------
# Edit 1: split "npm install async" into "npm install" / " --save" /
# " async", then drop "_GoBack" in between the last two runs.
# -----------------------------------------------------------------

$full = $d.Content.Text
$oldText = "npm install async"
$idx = $full.IndexOf($oldText)

$d.Range($idx, $idx + $oldText.Length).Text = "npm install"

$afterInstall = $idx + "npm install".Length
$d.Range($afterInstall, $afterInstall).InsertAfter(" --save")

# Protect the "npm install"/" --save" boundary so it doesn't get pulled
# back together when " async" is appended next.
$d.Bookmarks.Add("TEMP_PROTECT_B", $d.Range($afterInstall, $afterInstall))

$afterSave = $afterInstall + " --save".Length
$d.Range($afterSave, $afterSave).InsertAfter(" async")

$d.Bookmarks("TEMP_PROTECT_B").Delete()

# Finally, move "_GoBack" to its real new home: between " --save" and
# " async".
$d.Bookmarks.Add("_GoBack", $d.Range($afterSave, $afterSave))
